# Saving account Test cases added
#
# Renames "Sheet2" to "NewSavingAccount", populates it with two new
# saving-account test rows (CCC/4000, DDD/5000) mirroring the layout of
# the existing "NewCheckingAccount" sheet, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

$wsChecking = $wb.Worksheets.Item("NewCheckingAccount")
$wsSaving   = $wb.Worksheets.Item("Sheet2")

# --- Rename Sheet2 -> NewSavingAccount ------------------------------------
$wsSaving.Name = "NewSavingAccount"

# --- Write the new test data ----------------------------------------------
$wsSaving.Range("A1").Value = "Account Name"
$wsSaving.Range("B1").Value = "Deposit Amount"
$wsSaving.Range("A2").Value = "CCC"
$wsSaving.Range("B2").Value = 4000
$wsSaving.Range("A3").Value = "DDD"
$wsSaving.Range("B3").Value = 5000

# --- Match the look & feel of the checking-account sheet -------------------
# (header fill/border style + bordered data rows)
$wsChecking.Range("A1:B3").Copy()
$wsSaving.Range("A1:B3").PasteSpecial(-4122)  # xlPasteFormats

# Column widths to fit the new header/content text (same as the other sheet)
$wsSaving.Columns.Item(1).ColumnWidth = 13.166666666666666
$wsSaving.Columns.Item(2).ColumnWidth = 14.7

# --- Selection bookkeeping ---------------------------------------------
# Checking-account sheet keeps its data selected once it is no longer
# the active tab ...
$null = $wsChecking.Range("A1:B3").Select()

# ... and the new saving-account sheet becomes the active tab, with the
# cursor parked just past the data (C3).
$wsSaving.Activate()
$null = $wsSaving.Range("C3").Select()
